$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 599.6667
$ws.Range("I12").Value = 466
$ws.Range("J12").Value = 733.3333
$ws.Range("K12").Value = 466
$ws.Range("L12").Value = 733.3333
$ws.Range("M12").Value = -296
$ws.Range("N12").Value = -1073.3333

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()

$ws.Range("H38").Value = 4396.3335
$ws.Range("J38").Value = 6494.5
$ws.Range("L38").Value = 19483.5
$ws.Range("N38").Value = -20227.5

$ws.Range("H43").Value = 7950
$ws.Range("I43").Value = 7950
$ws.Range("K43").Value = 7950
$ws.Range("M43").Value = -7881

$ws.Range("H53").Value = 174.57143
$ws.Range("J53").Value = 427.5
$ws.Range("L53").Value = 427.5
$ws.Range("N53").Value = -1701.5

$ws.Range("H88").Value = 1699
$ws.Range("I88").Value = 1699
$ws.Range("K88").Value = 1699
$ws.Range("M88").Value = -1293

$ws.Range("H91").Value = 1699
$ws.Range("I91").Value = 1699
$ws.Range("K91").Value = 1699
$ws.Range("M91").Value = -295

$ws.Range("H98").Value = 874.125
$ws.Range("I98").Value = 874.125
$ws.Range("K98").Value = 874.125
$ws.Range("M98").Value = 623.875

$ws.Range("H107").Value = 976.625
$ws.Range("I107").Value = 941.73334
$ws.Range("K107").Value = 941.73334
$ws.Range("M107").Value = 978.26666

$ws.Range("H116").Value = 5058.9287
$ws.Range("I116").Value = 4992.6
$ws.Range("K116").Value = 4992.6
$ws.Range("M116").Value = -1550.6

$ws.Range("H122").Value = 874.125
$ws.Range("I122").Value = 874.125
$ws.Range("K122").Value = 2622.375
$ws.Range("M122").Value = -172.375

$ws.Range("H137").Value = 3999.75
$ws.Range("I137").Value = 999
$ws.Range("J137").Value = 5000
$ws.Range("K137").Value = 2997
$ws.Range("L137").Value = 15000
$ws.Range("M137").Value = -447
$ws.Range("N137").Value = -20100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3470.9583
$ws.Range("I32").Value = 3447.9565
$ws.Range("K32").Value = 3447.9565
$ws.Range("M32").Value = -3160.9565

$ws.Range("H74").Value = 2234.7646
$ws.Range("I74").Value = 666
$ws.Range("K74").Value = 666
$ws.Range("M74").Value = 208

$ws.Range("H77").Value = 2234.7646
$ws.Range("I77").Value = 666
$ws.Range("K77").Value = 3330
$ws.Range("M77").Value = 1038

$ws.Range("H122").Value = 1599.6364
$ws.Range("I122").Value = 1599.6364
$ws.Range("K122").Value = 4798.9092
$ws.Range("M122").Value = -2348.9092

$ws.Range("H132").Value = 2501.5
$ws.Range("I132").Value = 2078.5715
$ws.Range("K132").Value = 6235.7145
$ws.Range("M132").Value = -3705.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4065.6924
$ws.Range("I134").Value = 3905.818
$ws.Range("K134").Value = 11717.454
$ws.Range("M134").Value = -9182.454000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7254.6875
$ws.Range("I31").Value = 2118.8
$ws.Range("K31").Value = 2118.8
$ws.Range("M31").Value = -1823.8

$ws.Range("H34").Value = 7254.6875
$ws.Range("I34").Value = 2118.8
$ws.Range("K34").Value = 2118.8
$ws.Range("M34").Value = -1916.8

$ws.Range("H58").Value = 2048.8333
$ws.Range("I58").Value = 1878.8
$ws.Range("K58").Value = 1878.8
$ws.Range("M58").Value = -1675.8

$ws.Range("H107").Value = 967.75
$ws.Range("I107").Value = 728.9
$ws.Range("J107").Value = 2162
$ws.Range("K107").Value = 728.9
$ws.Range("L107").Value = 2162
$ws.Range("M107").Value = 1191.1
$ws.Range("N107").Value = -6002

$ws.Range("H122").Value = 1457.4286
$ws.Range("J122").Value = 1375
$ws.Range("L122").Value = 4125
$ws.Range("N122").Value = -9025

$ws.Range("H132").Value = 2209
$ws.Range("I132").Value = 1060.6471
$ws.Range("K132").Value = 3181.9413
$ws.Range("M132").Value = -651.9412999999995

$ws.Range("H133").Value = 124900
$ws.Range("J133").Value = 124900
$ws.Range("L133").Value = 124900
$ws.Range("N133").Value = -129960

$ws.Range("H136").Value = 2048.8333
$ws.Range("I136").Value = 1878.8
$ws.Range("K136").Value = 5636.4
$ws.Range("M136").Value = -3086.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1278.7142
$ws.Range("J68").Value = 1291.8334
$ws.Range("L68").Value = 3875.5002
$ws.Range("N68").Value = -5497.5002

$ws.Range("H71").Value = 1278.7142
$ws.Range("J71").Value = 1291.8334
$ws.Range("L71").Value = 11626.5006
$ws.Range("N71").Value = -19738.5006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6149.6924
$ws.Range("I80").Value = 5823.8335
$ws.Range("J80").Value = 6429
$ws.Range("K80").Value = 5823.8335
$ws.Range("L80").Value = 6429
$ws.Range("M80").Value = -4825.8335
$ws.Range("N80").Value = -8425

$ws.Range("H83").Value = 6149.6924
$ws.Range("I83").Value = 5823.8335
$ws.Range("J83").Value = 6429
$ws.Range("K83").Value = 29119.1675
$ws.Range("L83").Value = 32145
$ws.Range("M83").Value = -24127.1675
$ws.Range("N83").Value = -42129

$ws.Range("H92").Value = 19999.5
$ws.Range("J92").Value = 19999.5
$ws.Range("L92").Value = 19999.5
$ws.Range("N92").Value = -23743.5

$ws.Range("H95").Value = 48750
$ws.Range("J95").Value = 48750
$ws.Range("L95").Value = 48750
$ws.Range("N95").Value = -54242

$ws.Range("H102").Value = 4125.6665
$ws.Range("I102").Value = 3126.3333
$ws.Range("J102").Value = 5125
$ws.Range("K102").Value = 3126.3333
$ws.Range("L102").Value = 5125
$ws.Range("M102").Value = -1504.3333
$ws.Range("N102").Value = -8369

$ws.Range("H132").Value = 2951.2122
$ws.Range("I132").Value = 2375.36
$ws.Range("K132").Value = 7126.08
$ws.Range("M132").Value = -4596.08

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2999
$ws.Range("I46").Value = 2999
$ws.Range("K46").Value = 2999
$ws.Range("M46").Value = -2811

$ws.Range("H55").Value = 662.2222
$ws.Range("I55").Value = 581.4286
$ws.Range("K55").Value = 581.4286
$ws.Range("M55").Value = -408.4286

$ws.Range("H82").Value = 1798.3334
$ws.Range("I82").Value = 1700
$ws.Range("J82").Value = 1995
$ws.Range("K82").Value = 1700
$ws.Range("L82").Value = 1995
$ws.Range("M82").Value = -1339
$ws.Range("N82").Value = -2717

$ws.Range("H85").Value = 1798.3334
$ws.Range("I85").Value = 1700
$ws.Range("J85").Value = 1995
$ws.Range("K85").Value = 1700
$ws.Range("L85").Value = 1995
$ws.Range("M85").Value = -452
$ws.Range("N85").Value = -4491

$ws.Range("H101").Value = 20000
$ws.Range("J101").Value = 20000
$ws.Range("L101").Value = 20000
$ws.Range("N101").Value = -26490

$ws.Range("H122").Value = 2449
$ws.Range("I122").Value = 2449
$ws.Range("K122").Value = 7347
$ws.Range("M122").Value = -4897

$ws.Range("H132").Value = 3041.7368
$ws.Range("I132").Value = 2307.2856
$ws.Range("K132").Value = 6921.8568
$ws.Range("M132").Value = -4391.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1460.8334
$ws.Range("I122").Value = 1337.25
$ws.Range("J122").Value = 2449.5
$ws.Range("K122").Value = 4011.75
$ws.Range("L122").Value = 7348.5
$ws.Range("M122").Value = -1561.75
$ws.Range("N122").Value = -12248.5

$ws.Range("H132").Value = 1981.0238
$ws.Range("I132").Value = 1810.3549
$ws.Range("J132").Value = 2462
$ws.Range("K132").Value = 5431.0647
$ws.Range("L132").Value = 7386
$ws.Range("M132").Value = -2901.0647
$ws.Range("N132").Value = -12446

$ws.Range("H136").Value = 6256.6665
$ws.Range("I136").Value = 7766.5
$ws.Range("J136").Value = 1727.1666
$ws.Range("K136").Value = 23299.5
$ws.Range("L136").Value = 5181.4998
$ws.Range("M136").Value = -20749.5
$ws.Range("N136").Value = -10281.4998
